$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.328.72'
$ws.Range("E2").Value = '  -0.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.845.95'
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9978'
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.14'
$ws.Range("E5").Value = '  -0.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6262'
$ws.Range("E6").Value = '  -0.58%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9987'
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07586'
$ws.Range("E8").Value = '  -1.32%  '

$ws.Range("E9").Value = '  -1.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.64'
$ws.Range("E10").Value = '  +0.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07736'
$ws.Range("E11").Value = '  -0.16%  '

$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6778'
$ws.Range("E13").Value = '  -0.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001059'
$ws.Range("E14").Value = '  -2.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.84'
$ws.Range("E15").Value = '  -1.11%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.115'
$ws.Range("E16").Value = '  -0.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.359.76'
$ws.Range("E17").Value = '  -0.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.28'
$ws.Range("E18").Value = '  -1.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.33'
$ws.Range("E19").Value = '  -1.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9985'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.474'
$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9987'
$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("E24").Value = '  -0.37%  '

$ws.Range("E25").Value = '  +0.41%  '

$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.432'
$ws.Range("E27").Value = '  +8.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.459'
$ws.Range("E28").Value = '  -0.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05604'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.097'
$ws.Range("E30").Value = '  -0.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.056'
$ws.Range("E31").Value = '  +0.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.159'
$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("E33").Value = '  -1.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.6950'
$ws.Range("E34").Value = '  -1.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.581'
$ws.Range("E35").Value = '  -0.28%  '

$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.226.54'
$ws.Range("E36").Value = '  +0.14%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01794'
$ws.Range("E37").Value = '  -0.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.720'
$ws.Range("E38").Value = '  -2.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.353'
$ws.Range("E39").Value = '  -1.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8974'
$ws.Range("E40").Value = '  -1.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9983'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.34'
$ws.Range("E42").Value = '  -0.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.40'

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.190'
$ws.Range("E44").Value = '  +0.67%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000120'
$ws.Range("E45").Value = '  -0.97%  '

$ws.Range("E46").Value = '  -1.03%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.025'
$ws.Range("E47").Value = '  -0.19%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.686'
$ws.Range("E48").Value = '  -0.28%  '

$ws.Range("E49").Value = '  +1.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05693'
$ws.Range("E50").Value = '  -0.37%  '

$ws.Range("E51").Value = '  -0.17%  '
